# Resident 1.0.11 cases.xlsx -- "Add files via upload" edit
#
# Adds 10 new testcase rows (51-60) covering Resident service OTP,
# Policy details and Credential Issuance scenarios, plus the wrap-text
# styles those rows use, and refreshes the sheet view / window metadata
# to match the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Workbook-level bookkeeping that accompanies the authoring session
# (path of the machine that last saved the file, and the window size
# Excel remembered on close)
# ---------------------------------------------------------------------
try { $wb.Path = "C:\Users\Rakshit.B\Downloads\" } catch {}

$win = $excel.ActiveWindow
try { $win.Width = 20490 } catch {}
try { $win.Height = 7650 } catch {}

# ---------------------------------------------------------------------
# Row 51 -- Resident service_OTP_01
# ---------------------------------------------------------------------
$ws.Range("A51").Value = "Resident service_OTP_01"
$ws.Range("B51").Value = "Resident service"
$ws.Range("C51").Value = "OTP"
$ws.Range("D51").Value = "Functional"
$ws.Range("E51").Value = "Verify request OTP through new resident  OTP service"
$ws.Range("F51").Value = "Resident should receive OTP"
$ws.Range("G51").Value = "N"
$ws.Range("E51:F51").WrapText = $true
$ws.Rows.Item(51).RowHeight = 45

# ---------------------------------------------------------------------
# Row 52 -- Resident service_OTP_02
# ---------------------------------------------------------------------
$ws.Range("A52").Value = "Resident service_OTP_02"
$ws.Range("B52").Value = "Resident service"
$ws.Range("C52").Value = "OTP"
$ws.Range("D52").Value = "Functional"
$ws.Range("E52").Value = "Verify request OTP should expire after cofigured time"
$ws.Range("F52").Value = "Resident should receive OTP"
$ws.Range("G52").Value = "N"
$ws.Range("E52:F52").WrapText = $true
$ws.Rows.Item(52).RowHeight = 45

# ---------------------------------------------------------------------
# Row 53 -- Resident service_OTP_03
# ---------------------------------------------------------------------
$ws.Range("A53").Value = "Resident service_OTP_03"
$ws.Range("B53").Value = "Resident service"
$ws.Range("C53").Value = "OTP"
$ws.Range("D53").Value = "Functional"
$ws.Range("E53").Value = "Verify request OTP  using valid data"
$ws.Range("F53").Value = "Resident should receive OTP"
$ws.Range("G53").Value = "N"
$ws.Range("E53:F53").WrapText = $true
$ws.Rows.Item(53).RowHeight = 30

# ---------------------------------------------------------------------
# Row 54 -- Resident service_OTP_04
# ---------------------------------------------------------------------
$ws.Range("A54").Value = "Resident service_OTP_04"
$ws.Range("B54").Value = "Resident service"
$ws.Range("C54").Value = "OTP"
$ws.Range("D54").Value = "Functional"
$ws.Range("E54").Value = "Verify request OTP  using invalid data"
$ws.Range("F54").Value = "OTP should not be send to resident and should get appropriate error/ validation message"
$ws.Range("G54").Value = "N"
$ws.Range("E54:F54").WrapText = $true
$ws.Rows.Item(54).RowHeight = 45

# ---------------------------------------------------------------------
# Row 55 -- Resident service_OTP_05
# ---------------------------------------------------------------------
$ws.Range("A55").Value = "Resident service_OTP_05"
$ws.Range("B55").Value = "Resident service"
$ws.Range("C55").Value = "OTP"
$ws.Range("D55").Value = "Functional"
$ws.Range("E55").Value = "Verify OTP received  through EMAIL/MOBILE number"
$ws.Range("F55").Value = "Resident should receive OTP"
$ws.Range("G55").Value = "N"
$ws.Range("E55:F55").WrapText = $true
$ws.Rows.Item(55).RowHeight = 45

# ---------------------------------------------------------------------
# Row 56 -- Resident service_Policy details_01
# ---------------------------------------------------------------------
$ws.Range("A56").Value = "Resident service_Policy details_01"
$ws.Range("B56").Value = "Resident service"
$ws.Range("C56").Value = "Policy details"
$ws.Range("D56").Value = "Functional"
$ws.Range("E56").Value = "Mosip partner should be able to  fetch policy details using credential type and partner id"
$ws.Range("F56").Value = "Partner should get polcy details"
$ws.Range("G56").Value = "N"
$ws.Range("E56:F56").WrapText = $true
$ws.Rows.Item(56).RowHeight = 60

# ---------------------------------------------------------------------
# Row 57 -- Resident service_Policy details_02
# ---------------------------------------------------------------------
$ws.Range("A57").Value = "Resident service_Policy details_02"
$ws.Range("B57").Value = "Resident service"
$ws.Range("C57").Value = "Policy details"
$ws.Range("D57").Value = "Functional"
$ws.Range("E57").Value = "Verify by not giving mandatory input parameters(Credential Type and Partner ID)"
$ws.Range("F57").Value = "Partner should not get policy details"
$ws.Range("G57").Value = "N"
$ws.Range("E57:F57").WrapText = $true
$ws.Rows.Item(57).RowHeight = 60

# ---------------------------------------------------------------------
# Row 58 -- Resident service_Policy details_03
# ---------------------------------------------------------------------
$ws.Range("A58").Value = "Resident service_Policy details_03"
$ws.Range("B58").Value = "Resident service"
$ws.Range("C58").Value = "Policy details"
$ws.Range("D58").Value = "Functional"
$ws.Range("E58").Value = "Verify by giving invalid data in credential type"
$ws.Range("F58").Value = "Partner should not get policy details"
$ws.Range("G58").Value = "N"
$ws.Range("E58:F58").WrapText = $true
$ws.Rows.Item(58).RowHeight = 30

# ---------------------------------------------------------------------
# Row 59 -- Resident service_Policy details_04
# ---------------------------------------------------------------------
$ws.Range("A59").Value = "Resident service_Policy details_04"
$ws.Range("B59").Value = "Resident service"
$ws.Range("C59").Value = "Policy details"
$ws.Range("D59").Value = "Functional"
$ws.Range("E59").Value = "Verify by giving invalid data in partner id"
$ws.Range("F59").Value = "Partner should not get policy details"
$ws.Range("G59").Value = "N"
$ws.Range("E59:F59").WrapText = $true
$ws.Rows.Item(59).RowHeight = 30

# ---------------------------------------------------------------------
# Row 60 -- Resident Service_Credentialissuance_01
# ---------------------------------------------------------------------
$ws.Range("A60").Value = "Resident Service_Credentialissuance_01"
$ws.Range("B60").Value = "Resident Services"
$ws.Range("C60").Value = "`nCredential Issuance`n"
$ws.Range("D60").Value = "Functional"
$ws.Range("E60").Value = "Verify the print service with QR code having the best finger"
$ws.Range("F60").Value = "QR code should have best two fingers"
$ws.Range("G60").Value = "N"
$ws.Range("C60").WrapText = $true
$ws.Range("C60").HorizontalAlignment = -4131   # xlLeft
$ws.Range("C60").VerticalAlignment = -4160     # xlTop
$ws.Range("E60:F60").WrapText = $true
$ws.Rows.Item(60).RowHeight = 60

# ---------------------------------------------------------------------
# Refresh the view: new bottom row selected/visible
# ---------------------------------------------------------------------
[void]$ws.Range("H59").Select()
